$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("G22").Value = "wait(2);`nvalidate1;`nlink_Click(filemanagement_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT056_1041_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nwait(100);`nvalidate4;"
$ws.Range("G23").Value = "wait(2);`nvalidate1;`nlink_Click(filemanagement_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT056_1042_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nwait(100);`nvalidate4;"
$ws.Range("G24").Value = "wait(2);`nvalidate1;`nlink_Click(filemanagement_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT056_1043_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nwait(100);`nvalidate4;"
$ws.Range("G25").Value = "wait(2);`nvalidate1;`nlink_Click(filemanagement_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT056_1044_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nwait(100);`nvalidate4;"
